$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.985.03"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3
$ws.Range("D3").Value = "3.891.34"
$ws.Range("E3").Value = "  -2.21%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.37"
$ws.Range("E5").Value = "  +1.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.65"
$ws.Range("E6").Value = "  +10.60%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.669"
$ws.Range("E7").Value = "  -1.33%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.769"
$ws.Range("E9").Value = "  +2.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.180"
$ws.Range("E10").Value = "  +7.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.59"
$ws.Range("E11").Value = "  +1.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000324"
$ws.Range("E12").Value = "  +2.04%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.35"
$ws.Range("E13").Value = "  +4.51%  "

# Row 14
$ws.Range("D14").Value = "4.508.30"
$ws.Range("E14").Value = "  -1.94%  "

# Row 15
$ws.Range("D15").Value = "3.902.20"
$ws.Range("E15").Value = "  -1.83%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.16"
$ws.Range("E16").Value = "  +3.42%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.94"
$ws.Range("E17").Value = "  -0.61%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.21"
$ws.Range("E18").Value = "  -4.72%  "

# Row 19
$ws.Range("E19").Value = "  -1.59%  "

# Row 20
$ws.Range("D20").Value = "71.783.87"
$ws.Range("E20").Value = "  -0.83%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "437.69"
$ws.Range("E21").Value = "  +1.40%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.72"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.05"
$ws.Range("E23").Value = "  -1.75%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.29"
$ws.Range("E24").Value = "  -5.29%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.85"
$ws.Range("E25").Value = "  -2.75%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.23"
$ws.Range("E26").Value = "  -5.36%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.04"
$ws.Range("E27").Value = "  -3.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.93"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.21"
$ws.Range("E29").Value = "  -3.86%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.16"
$ws.Range("E30").Value = "  -3.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.85"
$ws.Range("E31").Value = "  -1.74%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.58"
$ws.Range("E32").Value = "  +0.58%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "49.85"
$ws.Range("E33").Value = "  +0.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.126"
$ws.Range("E34").Value = "  -4.59%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0984"
$ws.Range("E35").Value = "  +14.64%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.90"
$ws.Range("E36").Value = "  -2.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "617.10"
$ws.Range("E37").Value = "  -9.38%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.423"
$ws.Range("E38").Value = "  -3.80%  "

# Row 39
$ws.Range("E39").Value = "  -0.45%  "

# Row 40
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.34"
$ws.Range("E40").Value = "  +0.78%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.26%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.143"
$ws.Range("E42").Value = "  -2.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.20"
$ws.Range("E43").Value = "  +42.79%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0470"
$ws.Range("E44").Value = "  -3.39%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.23"
$ws.Range("E45").Value = "  -7.70%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.65"
$ws.Range("E46").Value = "  -6.04%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.144"
$ws.Range("E47").Value = "  -3.09%  "

# Row 48
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.84"
$ws.Range("E48").Value = "  -16.18%  "

# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.33"
$ws.Range("E49").Value = "  +0.16%  "

# Row 50
$ws.Range("D50").Value = "2.859.34"
$ws.Range("E50").Value = "  +2.25%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000279"
$ws.Range("E51").Value = "  +3.03%  "
